$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J
# Copy formatting from the existing header cell (H1) so I1/J1 match the
# bold/bordered/centered header style used by the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-17
$data = @(
    @(6, 9),
    @(13, 13),
    @(1, 4),
    @(1, 4),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(10, 10),
    @(6, 9),
    @(5, 8),
    @(5, 6),
    @(3, 7),
    @(1, 3),
    @(1, 4),
    @(5, 5),
    @(3, 3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
